$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old text -> corrected text (status name fix)
$map = @{
    "bleu" = "noir"
    "pas de résultat ni de publication" = "pas de résultat postés ni publiés"
    "résultat et / ou publication posté" = "résultat postés ou publiés"
    "résultat et / ou publication posté dans les 12 mois" = "résultat postés ou publiés dans les 12 mois"
    "résultat et / ou publication posté dans les 36 mois" = "résultat postés ou publiés dans les 36 mois"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count
$startRow = $used.Row
$startCol = $used.Column

for ($r = 0; $r -lt $rows; $r++) {
    for ($c = 0; $c -lt $cols; $c++) {
        $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
        $val = $cell.Value2
        if ($null -ne $val -and $map.ContainsKey($val)) {
            $cell.Value2 = $map[$val]
        }
    }
}
